$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("description (รายละเอียด)")

# Update the header text on the "description" sheet to include the extra
# clarifying note about entering numbers only / allowing blanks.
$ws2.Range("A1").Value2 = "เลขที่  (null) กรอกเป็นตัวเลขเท่านั้น สามารถเว้นว่างได้ "

# Change the selected cell on Sheet1 (it loses the "active/selected tab" flag)
$ws1.Range("D11").Select()

# Change the selected cell on the description sheet and make it the active
# (selected) tab of the workbook.
$ws2.Range("B5").Select()
$ws2.Activate()
